# Fixed data dictionary for users
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the table name from "User" to "users"
$ws.Range("A4").Value = "users"

# Update the title in A1, center it, and merge it across the table width
$ws.Range("A1").Value = "Data Dictionary for users"
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1:H1").MergeCells = $true
